# AP120_TestData_CreateAccountinginPayables_21C.xlsx
# Commit: "Add files via upload" / "Anu - AP Files Uploaded"
#
# The only substantive content change is on the "Input_Value" sheet:
# the credential values in R2:T2 (URL / UserName / Password) are
# cleared out before the file is (re-)uploaded. Clearing the cells
# (rather than overwriting with blanks) drops the now-unused shared
# strings from the workbook automatically, which is what shrinks the
# shared string table and shifts the other sheets' string indices.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

[void]$ws.Range("R2:T2").Select()
[void]$ws.Range("R2:T2").ClearContents()
